# Helper: convert an EMU value to the "points" value to feed into the COM
# Left/Top/Width/Height setters so that, after the host's float32 round-trip
# (which truncates rather than rounds when converting back to EMU), the
# saved OOXML ends up with exactly the intended EMU value.
function EMU($emu) {
    return ($emu + 0.5) / 12700
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1. The existing "Oval 20" shape is duplicated to create a brand new,
#    larger oval ("Oval 24") with identical styling, which is then moved
#    to its own target location and sent behind everything else.
# ---------------------------------------------------------------------
$oval20 = $s.Shapes.Item("Oval 20")

$oval24 = $oval20.Duplicate().Item(1)
$oval24.Name = "Oval 24"
$oval24.Left = EMU(7206742)
$oval24.Top = EMU(1699099)
$oval24.Width = EMU(1897496)
$oval24.Height = EMU(1309990)

# ---------------------------------------------------------------------
# 2. The original "Oval 20" shape shifts slightly to the left and is
#    also sent near the back (just above the new "Oval 24").
# ---------------------------------------------------------------------
$oval20.Left = EMU(7554472)

# z-order: send Oval 20 to the back first, then Oval 24, so that the
# final back-to-front order is: Oval 24, Oval 20, <everything else>.
$oval20.ZOrder(1)   # msoSendToBack
$oval24.ZOrder(1)   # msoSendToBack

# ---------------------------------------------------------------------
# 3. "TextBox 15" (currently reading "outer") moves slightly and grows
#    a bit wider, and its text is changed to "atomic".
# ---------------------------------------------------------------------
$tb15 = $s.Shapes.Item("TextBox 15")
$tb15.Left = EMU(7633541)
$tb15.Width = EMU(539315)
$tb15.TextFrame.TextRange.Text = "atomic"

# ---------------------------------------------------------------------
# 4. "Oval 16" shifts to the right.
# ---------------------------------------------------------------------
$oval16 = $s.Shapes.Item("Oval 16")
$oval16.Left = EMU(8516535)

# ---------------------------------------------------------------------
# 5. "TextBox 17" (reading "absent") shifts to the right.
# ---------------------------------------------------------------------
$tb17 = $s.Shapes.Item("TextBox 17")
$tb17.Left = EMU(9161077)

# ---------------------------------------------------------------------
# 6. A new label "outer" is introduced (taking over the role vacated by
#    the now-relocated "Oval 20"), cloned from "TextBox 15" so it keeps
#    the same italic / colour / no-fill styling, then repositioned,
#    retitled and retexted.
# ---------------------------------------------------------------------
$newOuter = $tb15.Duplicate().Item(1)
$newOuter.Name = "TextBox 25"
$newOuter.Left = EMU(7951556)
$newOuter.Top = EMU(1723400)
$newOuter.Width = EMU(407869)
$newOuter.Height = EMU(215444)
$newOuter.TextFrame.TextRange.Text = "outer"
